$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.760.54"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "3.409.59"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.35%  "

# Row 5
$ws.Range("D5").Value = "'412.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "

# Row 6
$ws.Range("D6").Value = "'130.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "

# Row 7
$ws.Range("D7").Value = "'0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.38%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.725"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").Value = "'0.134"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.84%  "

# Row 11
$ws.Range("D11").Value = "'42.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.02%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'9.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.69%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.945.86"
$ws.Range("E13").Value = "  -0.05%  "

# Row 14
$ws.Range("E14").Value = "  +0.10%  "

# Row 15
$ws.Range("D15").Value = "'0.0000209"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.58%  "

# Row 16
$ws.Range("D16").Value = "'20.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.74%  "

# Row 17
$ws.Range("D17").Value = "3.413.10"
$ws.Range("E17").Value = "  +0.19%  "

# Row 18
$ws.Range("E18").Value = "  +1.76%  "

# Row 19
$ws.Range("D19").Value = "'12.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.07%  "

# Row 20
$ws.Range("D20").Value = "61.757.27"
$ws.Range("E20").Value = "  -0.12%  "

# Row 21
$ws.Range("D21").Value = "'478.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +18.25%  "

# Row 22
$ws.Range("D22").Value = "'89.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "

# Row 23
$ws.Range("D23").Value = "'3.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.22%  "

# Row 24
$ws.Range("D24").Value = "'13.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "

# Row 25
$ws.Range("D25").Value = "'3.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.44%  "

# Row 26
$ws.Range("D26").Value = "'9.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.30%  "

# Row 27
$ws.Range("D27").Value = "'33.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.11%  "

# Row 28
$ws.Range("D28").Value = "'4.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "

# Row 29
$ws.Range("D29").Value = "'7.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.20%  "

# Row 30
$ws.Range("D30").Value = "'11.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.77%  "

# Row 31
$ws.Range("D31").Value = "'2.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.31%  "

# Row 32
$ws.Range("E32").Value = "  -1.68%  "

# Row 33
$ws.Range("E33").Value = "  -4.00%  "

# Row 34
$ws.Range("D34").Value = "'40.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.68%  "

# Row 35
$ws.Range("E35").Value = "  -0.76%  "

# Row 36
$ws.Range("D36").Value = "'57.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.12%  "

# Row 37
$ws.Range("D37").Value = "'0.0486"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.37%  "

# Row 38
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "

# Row 39
$ws.Range("E39").Value = "  +4.35%  "

# Row 40
$ws.Range("D40").Value = "'0.326"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.33%  "

# Row 41
$ws.Range("D41").Value = "'147.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.05%  "

# Row 42
$ws.Range("E42").Value = "  +1.15%  "

# Row 43
$ws.Range("E43").Value = "  -0.16%  "

# Row 44
$ws.Range("D44").Value = "'2.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.27%  "

# Row 45
$ws.Range("E45").Value = "  +8.16%  "

# Row 46
$ws.Range("E46").Value = "  +4.56%  "

# Row 47
$ws.Range("D47").Value = "'2.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +19.12%  "

# Row 48
$ws.Range("D48").Value = "'16.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "

# Row 49
$ws.Range("D49").Value = "'22.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.46%  "

# Row 50
$ws.Range("E50").Value = "  +8.41%  "

# Row 51
$ws.Range("D51").Value = "'112.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.50%  "
